$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45609.5
$ws.Range("J3").Value = 45609.5
$ws.Range("L3").Value = 45609.5
$ws.Range("N3").Value = -45837.5
$ws.Range("H51").Value = 2516.6667
$ws.Range("I51").Value = 1766.6666
$ws.Range("K51").Value = 1766.6666
$ws.Range("M51").Value = -1282.6666
$ws.Range("H102").Value = 45609.5
$ws.Range("J102").Value = 45609.5
$ws.Range("L102").Value = 45609.5
$ws.Range("N102").Value = -52099.5
$ws.Range("H132").Value = 4236.476
$ws.Range("I132").Value = 4087.7368
$ws.Range("J132").Value = 5649.5
$ws.Range("K132").Value = 12263.2104
$ws.Range("L132").Value = 16948.5
$ws.Range("M132").Value = -9733.2104
$ws.Range("N132").Value = -22008.5
$ws.Range("H138").Value = 3570.7222
$ws.Range("I138").Value = 3001.3333
$ws.Range("J138").Value = 3789.718
$ws.Range("K138").Value = 9003.999899999999
$ws.Range("L138").Value = 11369.154
$ws.Range("M138").Value = -3863.999899999999
$ws.Range("N138").Value = -21649.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 86676.664
$ws.Range("J29").Value = 86676.664
$ws.Range("L29").Value = 86676.664
$ws.Range("N29").Value = -87292.664
$ws.Range("H31").Value = 20244.1
$ws.Range("I31").Value = 4634.4287
$ws.Range("K31").Value = 4634.4287
$ws.Range("M31").Value = -4340.4287
$ws.Range("H45").Value = 2030.9
$ws.Range("I45").Value = 1783.25
$ws.Range("J45").Value = 2402.375
$ws.Range("K45").Value = 1783.25
$ws.Range("L45").Value = 2402.375
$ws.Range("M45").Value = -1406.25
$ws.Range("N45").Value = -3156.375
$ws.Range("H61").Value = 3054.9607
$ws.Range("I61").Value = 1636.1724
$ws.Range("J61").Value = 4925.1816
$ws.Range("K61").Value = 1636.1724
$ws.Range("L61").Value = 4925.1816
$ws.Range("M61").Value = -1424.1724
$ws.Range("N61").Value = -5349.1816
$ws.Range("H122").Value = 103043.9
$ws.Range("I122").Value = 251409.75
$ws.Range("J122").Value = 4133.3335
$ws.Range("K122").Value = 754229.25
$ws.Range("L122").Value = 12400.0005
$ws.Range("M122").Value = -751779.25
$ws.Range("N122").Value = -17300.0005
$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229
$ws.Range("H132").Value = 1712682.9
$ws.Range("I132").Value = 2672.6296
$ws.Range("J132").Value = 4277698.5
$ws.Range("K132").Value = 8017.888800000001
$ws.Range("L132").Value = 12833095.5
$ws.Range("M132").Value = -5487.888800000001
$ws.Range("N132").Value = -12838155.5
$ws.Range("H136").Value = 3054.9607
$ws.Range("I136").Value = 1636.1724
$ws.Range("J136").Value = 4925.1816
$ws.Range("K136").Value = 4908.5172
$ws.Range("L136").Value = 14775.5448
$ws.Range("M136").Value = -2358.5172
$ws.Range("N136").Value = -19875.5448
$ws.Range("H139").Value = 63200.91
$ws.Range("J139").Value = 63200.91
$ws.Range("L139").Value = 63200.91
$ws.Range("N139").Value = -73480.91
$ws.Range("H141").Value = 87894.5
$ws.Range("J141").Value = 87894.5
$ws.Range("L141").Value = 87894.5
$ws.Range("N141").Value = -98254.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 3000
$ws.Range("J30").Value = 3000
$ws.Range("L30").Value = 3000
$ws.Range("N30").Value = -3250
$ws.Range("H92").Value = 92500
$ws.Range("J92").Value = 92500
$ws.Range("L92").Value = 92500
$ws.Range("N92").Value = -97492
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H97").Value = 9844.666999999999
$ws.Range("I97").Value = 2825.25
$ws.Range("J97").Value = 66000
$ws.Range("K97").Value = 2825.25
$ws.Range("L97").Value = 66000
$ws.Range("M97").Value = -1834.25
$ws.Range("N97").Value = -67982
$ws.Range("H102").Value = 28506.111
$ws.Range("I102").Value = 20194.375
$ws.Range("K102").Value = 20194.375
$ws.Range("M102").Value = -16949.375
$ws.Range("H134").Value = 2872
$ws.Range("I134").Value = 2681.0344
$ws.Range("J134").Value = 3564.25
$ws.Range("K134").Value = 8043.1032
$ws.Range("L134").Value = 10692.75
$ws.Range("M134").Value = -5508.1032
$ws.Range("N134").Value = -15762.75
$ws.Range("H138").Value = 10000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 70770
$ws.Range("J140").Value = 70770
$ws.Range("L140").Value = 70770
$ws.Range("N140").Value = -81130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1891.8572
$ws.Range("J16").Value = 1872.8
$ws.Range("L16").Value = 1872.8
$ws.Range("N16").Value = -2446.8
$ws.Range("H107").Value = 625.6539
$ws.Range("I107").Value = 333.41666
$ws.Range("J107").Value = 876.1429000000001
$ws.Range("K107").Value = 333.41666
$ws.Range("L107").Value = 876.1429000000001
$ws.Range("M107").Value = 1586.58334
$ws.Range("N107").Value = -4716.1429
$ws.Range("H113").Value = 1891.8572
$ws.Range("J113").Value = 1872.8
$ws.Range("L113").Value = 1872.8
$ws.Range("N113").Value = -6212.8
$ws.Range("H132").Value = 2033.6562
$ws.Range("I132").Value = 1714.6
$ws.Range("J132").Value = 2565.4167
$ws.Range("K132").Value = 5143.799999999999
$ws.Range("L132").Value = 7696.250100000001
$ws.Range("M132").Value = -2613.799999999999
$ws.Range("N132").Value = -12756.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4119664.8
$ws.Range("I4").Value = 17500632
$ws.Range("J4").Value = 2444.2307
$ws.Range("K4").Value = 52501896
$ws.Range("L4").Value = 7332.6921
$ws.Range("M4").Value = -52501784
$ws.Range("N4").Value = -7556.6921
$ws.Range("H12").Value = 95.90476
$ws.Range("I12").Value = 52.142857
$ws.Range("J12").Value = 117.78571
$ws.Range("K12").Value = 156.428571
$ws.Range("L12").Value = 353.35713
$ws.Range("M12").Value = 16.57142899999999
$ws.Range("N12").Value = -699.35713
$ws.Range("H49").Value = 6891
$ws.Range("J49").Value = 6891
$ws.Range("L49").Value = 20673
$ws.Range("N49").Value = -20985
$ws.Range("H68").Value = 1286.4615
$ws.Range("J68").Value = 1340.25
$ws.Range("L68").Value = 4020.75
$ws.Range("N68").Value = -5642.75
$ws.Range("H71").Value = 1286.4615
$ws.Range("J71").Value = 1340.25
$ws.Range("L71").Value = 12062.25
$ws.Range("N71").Value = -20174.25
$ws.Range("H115").Value = 6100.375
$ws.Range("I115").Value = 4003.3333
$ws.Range("J115").Value = 7358.6
$ws.Range("K115").Value = 12009.9999
$ws.Range("L115").Value = 22075.8
$ws.Range("M115").Value = -10834.9999
$ws.Range("N115").Value = -24425.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 54929.79
$ws.Range("I97").Value = 79110.38
$ws.Range("J97").Value = 2538.5
$ws.Range("K97").Value = 79110.38
$ws.Range("L97").Value = 2538.5
$ws.Range("M97").Value = -78614.38
$ws.Range("N97").Value = -3530.5
$ws.Range("H116").Value = 39999
$ws.Range("J116").Value = 39999
$ws.Range("L116").Value = 39999
$ws.Range("N116").Value = -49177
$ws.Range("H122").Value = 16785.715
$ws.Range("I122").Value = 27625
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 82875
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -80425
$ws.Range("N122").Value = -11899.9999
$ws.Range("H123").Value = 8324.190000000001
$ws.Range("J123").Value = 8324.190000000001
$ws.Range("L123").Value = 8324.190000000001
$ws.Range("N123").Value = -13224.19
$ws.Range("H132").Value = 2145.3845
$ws.Range("I132").Value = 1510.1818
$ws.Range("J132").Value = 2967.4119
$ws.Range("K132").Value = 4530.5454
$ws.Range("L132").Value = 8902.235700000001
$ws.Range("M132").Value = -2000.5454
$ws.Range("N132").Value = -13962.2357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4692.3335
$ws.Range("I7").Value = 4594.3335
$ws.Range("J7").Value = 5133.3335
$ws.Range("K7").Value = 4594.3335
$ws.Range("L7").Value = 5133.3335
$ws.Range("M7").Value = -4482.3335
$ws.Range("N7").Value = -5357.3335
$ws.Range("H14").Value = 82671.336
$ws.Range("J14").Value = 82671.336
$ws.Range("L14").Value = 82671.336
$ws.Range("N14").Value = -83015.336
$ws.Range("H61").Value = 3764.111
$ws.Range("J61").Value = 5500
$ws.Range("L61").Value = 5500
$ws.Range("N61").Value = -5904
$ws.Range("H113").Value = 3764.111
$ws.Range("J113").Value = 5500
$ws.Range("L113").Value = 5500
$ws.Range("N113").Value = -9840
$ws.Range("H122").Value = 4527
$ws.Range("I122").Value = 4777.778
$ws.Range("J122").Value = 4244.875
$ws.Range("K122").Value = 14333.334
$ws.Range("L122").Value = 12734.625
$ws.Range("M122").Value = -11883.334
$ws.Range("N122").Value = -17634.625
$ws.Range("H126").Value = 4692.3335
$ws.Range("I126").Value = 4594.3335
$ws.Range("J126").Value = 5133.3335
$ws.Range("K126").Value = 13783.0005
$ws.Range("L126").Value = 15400.0005
$ws.Range("M126").Value = -11313.0005
$ws.Range("N126").Value = -20340.0005
$ws.Range("H132").Value = 3052.611
$ws.Range("I132").Value = 2821.926
$ws.Range("J132").Value = 3744.6667
$ws.Range("K132").Value = 8465.778
$ws.Range("L132").Value = 11234.0001
$ws.Range("M132").Value = -5935.778
$ws.Range("N132").Value = -16294.0001
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1705.1666
$ws.Range("I113").Value = 2397.4285
$ws.Range("J113").Value = 736
$ws.Range("K113").Value = 7192.2855
$ws.Range("L113").Value = 2208
$ws.Range("M113").Value = -5022.2855
$ws.Range("N113").Value = -6548
$ws.Range("H132").Value = 6483871.5
$ws.Range("I132").Value = 2107.2424
$ws.Range("J132").Value = 24308722
$ws.Range("K132").Value = 6321.7272
$ws.Range("L132").Value = 72926166
$ws.Range("M132").Value = -3791.7272
$ws.Range("N132").Value = -72931226
